$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"1.141077787814027e-09"
$ws.Range("H2").Value = [double]"5.159229751235991e-08"
$ws.Range("K2").Value = 5.786861712847108
$ws.Range("L2").Value = "[3.7587098331745654, 7.81501359251965]"
$ws.Range("M2").Value = [double]"3.498788281319776e-08"
$ws.Range("N2").Value = [double]"6.997576562639551e-08"
$ws.Range("O2").Value = -1.38368445203154
$ws.Range("P2").Value = "[-1.8113687372049245, -0.9560001668581553]"
$ws.Range("Q2").Value = [double]"4.811149256767067e-10"
$ws.Range("R2").Value = [double]"4.811149256767067e-10"
$ws.Range("S2").Value = 10.57268868283663
$ws.Range("T2").Value = "[9.344280666675832, 11.801096698997437]"
$ws.Range("W2").Value = 5.71691691691705
$ws.Range("X2").Value = 3.949869869869963
$ws.Range("Y2").Value = 7.483963963964137

# Row 3
$ws.Range("E3").Value = 22.6700000000001
$ws.Range("G3").Value = 0.0002668698277130677
$ws.Range("H3").Value = 0.001057968345742211
$ws.Range("K3").Value = 6.003696598938154
$ws.Range("L3").Value = "[2.3393002808404173, 9.66809291703589]"
$ws.Range("M3").Value = 0.001450976671246051
$ws.Range("N3").Value = 0.001450976671246051
$ws.Range("O3").Value = 2.861711025792505
$ws.Range("P3").Value = "[2.1698687997767347, 3.5535532518082746]"
$ws.Range("Q3").Value = [double]"4.618527782440651e-14"
$ws.Range("R3").Value = [double]"9.237055564881302e-14"
$ws.Range("S3").Value = 11.30316861127206
$ws.Range("T3").Value = "[9.29235993288776, 13.313977289656357]"
$ws.Range("W3").Value = 12.34482482482488
$ws.Range("X3").Value = 9.848628628628672
$ws.Range("Y3").Value = 14.84102102102109
